$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$ws.Range("D2:D3").Value = -0.111
$ws.Range("E2:E3").Value = -0.07969999999999999
$ws.Range("G2:G3").Value = 0.9153846153846155
$ws.Range("H2:H3").Value = 0.9153846153846155
$ws.Range("I2:I3").Value = 0.6215384615384616
$ws.Range("J2:J3").Value = 0.5740090497737558
$ws.Range("K2:K3").Value = 6.58
$ws.Range("L2:L3").Value = 0.5061538461538462
$ws.Range("M2:M3").Value = 8.3607
$ws.Range("N2:N3").Value = 0.06802847843775427
$ws.Range("O2:O3").Value = 1.270623100303951
$ws.Range("P2:P3").Value = 8.3607
$ws.Range("Q2:Q3").Value = 0.06802847843775427
$ws.Range("R2:R3").Value = 1.270623100303951
$ws.Range("U2:U3").Value = 37
$ws.Range("V2:V3").Value = 0.3010577705451586
$ws.Range("W2:W3").Value = 0.03345195729537367
$ws.Range("X2:X3").Value = 0.03856114922038012
$ws.Range("Y2:Y3").Value = -0.005109191925006452
$ws.Range("Z2:Z3").Value = 0.0775499003782049
$ws.Range("AA2:AA3").Value = 0.04451434462614282
$ws.Range("AB2:AB3").Value = 0.03848216485452512
$ws.Range("AC2:AC3").Value = 0.006032179771617692
$ws.Range("AD2:AD3").Value = 0.515
$ws.Range("AE2:AE3").Value = 0
$ws.Range("AF2:AF3").Value = 0.515
$ws.Range("AG2:AG3").Value = -36.485
$ws.Range("AH2:AH3").Value = 0.004172912530891706
$ws.Range("AI2:AI3").Value = 0.002530526005454144
$ws.Range("AJ2:AJ3").Value = -0.4222067928021755
$ws.Range("AK2:AK3").Value = -0.2191093895444855
$ws.Range("AL2:AL3").Value = 0.006
$ws.Range("AM2:AM3").Value = 0.006
$ws.Range("AN2:AN3").Value = 0.06358024691358025
$ws.Range("AO2:AO3").Value = 1346.666666666667
$ws.Range("AP2:AP3").Value = -4.504320987654321
$ws.Range("AQ2:AQ3").Value = 1346.666666666667
